$d = $word.ActiveDocument

# Paragraph 1 is the hidden "**ID__...__ID**" bookmark-style marker paragraph
# at the very top of the document.
$p = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt clearance and no
# visible line - this matches the <w:pBdr> with only w:space attributes.
$borders = $p.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p.Format.LeftIndent = 11.25

# Replace the marker text and swallow the trailing single-space run that
# followed it, leaving one run with the updated marker text.
$range = $p.Range
$range.Find.Execute("**ID__AFFARS_pgi_5301_topic_4__ID** ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "**ID__AFFARS_AF_PGI_5301_170_2__ID**", 2)
